$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 changes from the text "R40" to the text "1".
# Prefix with an apostrophe so the numeric-looking value "1" is kept as
# text (a new shared-string entry) instead of being auto-converted to a
# number, matching the original cell's text type and style.
$ws.Range("B11").Value = "'1"
